# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    3 = @(1.445647641019636, 9.983522426115931, 3993.344853322108, 2797.565817734744, 6802.339841123987)
    4 = @(0.6545652718822623, 1.626987699542094, 189.6080260415259, 13.86384647080068, 205.753425483751)
    5 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.964545797025059)
    6 = @(0.6545652718822623, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.716211508195562)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G
}
